$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '60.349.14'
$ws.Range('E2').Value = '  +0.94%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.603.29'
$ws.Range('E3').Value = '  -0.03%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '576.11'
$ws.Range('E5').Value = '  +3.27%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '143.14'
$ws.Range('E6').Value = '  +0.85%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  +0.12%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.597'
$ws.Range('E8').Value = '  -0.47%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '2.608.29'
$ws.Range('E9').Value = '  -0.56%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '6.56'
$ws.Range('E10').Value = '  -2.04%  '
$ws.Range('E11').Value = '  +0.70%  '
$ws.Range('E12').Value = '  -3.70%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.369'
$ws.Range('E13').Value = '  +0.54%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '3.058.85'
$ws.Range('E14').Value = '  -0.02%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '24.40'
$ws.Range('E15').Value = '  +3.61%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '60.336.52'
$ws.Range('E16').Value = '  +0.99%  '
$ws.Range('E17').Value = '  +2.26%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '2.607.34'
$ws.Range('E18').Value = '  +0.13%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '11.43'
$ws.Range('E19').Value = '  +7.35%  '
$ws.Range('E20').Value = '  +0.35%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '346.90'
$ws.Range('E21').Value = '  +1.18%  '
$ws.Range('E22').Value = '  +1.54%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  -0.03%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '0.531'
$ws.Range('E24').Value = '  +1.45%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '63.17'
$ws.Range('E25').Value = '  +1.16%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.998'
$ws.Range('E26').Value = '  +0.12%  '
$ws.Range('E27').Value = '  +0.01%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '8.05'
$ws.Range('E28').Value = '  +6.05%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.0₃0796'
$ws.Range('E29').Value = '  +1.54%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.86'
$ws.Range('E30').Value = '  +9.73%  '
$ws.Range('E31').Value = '  +3.15%  '
$ws.Range('E32').Value = '  +0.06%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '166.40'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '19.43'
$ws.Range('E34').Value = '  +0.22%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '4.28'
$ws.Range('E35').Value = '  +3.59%  '
$ws.Range('E36').Value = '  +9.51%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.989'
$ws.Range('E37').Value = '  +7.41%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.62'
$ws.Range('E38').Value = '  +6.75%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '38.09'
$ws.Range('E39').Value = '  +0.91%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '313.13'
$ws.Range('E40').Value = '  +7.04%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '3.89'
$ws.Range('E41').Value = '  +4.74%  '
$ws.Range('E42').Value = '  -1.02%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '135.09'
$ws.Range('E43').Value = '  -3.16%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0994'
$ws.Range('E44').Value = '  +1.43%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.999'
$ws.Range('E45').Value = '  +0.24%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '19.85'
$ws.Range('E46').Value = '  +3.08%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0552'
$ws.Range('E47').Value = '  +2.27%  '
$ws.Range('E48').Value = '  +4.15%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.605'
$ws.Range('E49').Value = '  +0.47%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '20.01'
$ws.Range('E50').Value = '  +5.34%  '
$ws.Range('E51').Value = '  +0.44%  '
